$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 338. This shifts the existing rows 338..371
# down to 339..372, preserving all of their values (matching the diff,
# where the old row 338 data now appears at row 339, old row 339 at row
# 340, ... and old row 371 at row 372).
$ws.Rows(338).Insert()

# Populate the newly inserted (blank) row 338 with the same values as the
# row that used to occupy 338 (now at row 339), except for the date
# (column D) and volume (column J), which take on new values.
$ws.Cells.Item(338, 1).Value = 4
$ws.Cells.Item(338, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(338, 3).Value = "Los Lagos"
$ws.Cells.Item(338, 4).Value = 44769
$ws.Cells.Item(338, 5).Value = 10
$ws.Cells.Item(338, 6).Value = 100112023
$ws.Cells.Item(338, 7).Value = "Brócoli"
$ws.Cells.Item(338, 8).Value = "Sin especificar"
$ws.Cells.Item(338, 9).Value = "Primera"
$ws.Cells.Item(338, 10).Value = 250
$ws.Cells.Item(338, 11).Value = 1500
$ws.Cells.Item(338, 12).Value = 1500
$ws.Cells.Item(338, 13).Value = 1500
$ws.Cells.Item(338, 14).Value = "`$/unidad"
$ws.Cells.Item(338, 15).Value = "Región Metropolitana"
$ws.Cells.Item(338, 16).Value = 1500
$ws.Cells.Item(338, 17).Value = 1
$ws.Cells.Item(338, 18).Value = "Hortaliza"

# Match the date number format used by the rest of column D.
$ws.Cells.Item(338, 4).NumberFormat = $ws.Cells.Item(339, 4).NumberFormat
